$wb = $excel.ActiveWorkbook

# --- Sheet "validDatas" (sheet1) ---
$ws1 = $wb.Worksheets.Item(1)

# Row 4: replace old junk data with new junk data
$ws1.Range("A4").Value = "mndbmfnb"
$ws1.Range("B4").Value = "dsfbsdnfb"

# Row 5: brand-new row of data
$ws1.Range("A5").Value = "fjsdfh"
$ws1.Range("B5").Value = "hkjhkjh"

# Row 6: brand-new row of data
$ws1.Range("A6").Value = "yweiruy"
$ws1.Range("B6").Value = "sdhfkjhsd"

# --- Sheet "InvalidDatas" (sheet2) ---
$ws2 = $wb.Worksheets.Item(2)

# Select entire row 3 (no value changes on this sheet)
[void]$ws2.Rows.Item(3).Select()

# --- Re-activate sheet1 and set its selection to B6 ---
[void]$ws1.Activate()
[void]$ws1.Range("B6").Select()
